$p = $ppt.ActivePresentation
$newDate = "9/8/19 8:28 PM"

# Notes Master - the "Date Placeholder" shape that renders the cached
# datetime8 field shown on every notes page / notes master.
$nm = $p.NotesMaster
$nm.HeadersFooters.DateAndTime.Text = $newDate

# Handout Master - same date field, used when printing handouts.
$hm = $p.HandoutMaster
$hm.HeadersFooters.DateAndTime.Text = $newDate

# Each slide's Notes Page has its own copy of the cached field text
# inside its "Date Placeholder" shape (always shape index 5 in this deck).
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $ns = $p.Slides.Item($i).NotesPage
    $dateShape = $ns.Shapes.Item(5)
    $dateShape.TextFrame.TextRange.Text = $newDate
}
